$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FREQ")

# Column U "RMST" values (P=presente, F=falta) for rows 3-20, mirroring the
# grouping pattern used for column T (PROJ/PDMO/RMST per missao/grupo).
$valoresU = @{
    3  = "P"
    4  = "P"
    5  = "F"
    6  = "P"
    7  = "F"
    8  = "F"
    9  = "P"
    10 = "P"
    11 = "P"
    12 = "F"
    13 = "P"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "P"
    18 = "F"
    19 = "P"
    20 = "P"
}

foreach ($linha in $valoresU.Keys) {
    $ws.Range("U$linha").Value = $valoresU[$linha]
}

# Match column T's black-font formatting on the newly filled column U (and
# re-assert it on T, since both end up sharing the same style).
$ws.Range("T3:T20,U3:U20").Font.Color = 0

# Update the active selection to reflect where editing left off.
$ws.Range("U9").Select()

# Page layout: A4 portrait for printing.
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9
